# Update the "想去人数" (number of people interested) figures for three
# events that appear on both the "展览" sheet and the "全部类型" sheet.
$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 949
$wsExhibition.Range("F3").Value = 1831
$wsExhibition.Range("F4").Value = 414

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 949
$wsAll.Range("F5").Value = 1831
$wsAll.Range("F6").Value = 414
